$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new table row right after the current last data row (row 18) ---
# This pushes the existing row 18 down to row 19 (carrying its data/format),
# and shifts the lower block (old rows 19-24, blank gap + signature block) down by one.
$ws.Rows(19).Insert()

# The freshly inserted row 19 has default formatting; copy the *original* last-row
# formatting (still sitting on row 18) down onto it.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# Row 18 now becomes a "middle" table row like rows 16/17 - copy that formatting onto it.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Update the period labels shown for each table row (newest period added) ---
$ws.Range("E16").Value = "2507"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2505"

# --- Fill in the newly added table row with the (oldest) period's data ---
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73550995"
$ws.Range("D19").Value = "FARID ALBERTO ARROYO FERNANDEZ"
$ws.Range("E19").Value = "2504"
$ws.Range("F19").Value = 89540
$ws.Range("G19").Value = 2238500

# --- Update the summary totals: total overdue amount and number of periods ---
$ws.Range("E11").Value = 358160
$ws.Range("F13").Value = 4
